# Update the ProductLoanInput sheet: change "repaymentstrategy" value (B17)
# from "Mifos style" to "Penalties, Fees, Interest, Principal order", and
# update the active cell selection to B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"
$ws.Range("B1").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Activate()
$ws.Range("B17").Select()
